$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value = 17733.334
$ws.Range("I118").Value = 50650
$ws.Range("J118").Value = 1275
$ws.Range("K118").Value = 151950
$ws.Range("L118").Value = 3825
$ws.Range("M118").Value = -150293
$ws.Range("N118").Value = -7139

$ws.Range("H141").Value = 1769.7843
$ws.Range("I141").Value = 1609.9783
$ws.Range("J141").Value = 3240
$ws.Range("K141").Value = 4829.9349
$ws.Range("L141").Value = 9720
$ws.Range("M141").Value = 350.0650999999998
$ws.Range("N141").Value = -20080

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19409.494
$ws.Range("I32").Value = 4309.906
$ws.Range("J32").Value = 126364.914
$ws.Range("K32").Value = 4309.906
$ws.Range("L32").Value = 126364.914
$ws.Range("M32").Value = -4022.906
$ws.Range("N32").Value = -126938.914

$ws.Range("H41").Value = 11404.5
$ws.Range("I41").Value = 3852
$ws.Range("J41").Value = 34062
$ws.Range("K41").Value = 3852
$ws.Range("L41").Value = 34062
$ws.Range("M41").Value = -3438

$ws.Range("H61").Value = 2700.8064
$ws.Range("I61").Value = 1983.1666
$ws.Range("J61").Value = 2873.04
$ws.Range("K61").Value = 1983.1666
$ws.Range("L61").Value = 2873.04
$ws.Range("M61").Value = -1771.1666
$ws.Range("N61").Value = -3297.04

$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

$ws.Range("H122").Value = 1942.0385
$ws.Range("I122").Value = 1939.72
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 5819.16
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -3369.16
$ws.Range("N122").Value = -10900

$ws.Range("H132").Value = 2681.6667
$ws.Range("I132").Value = 2267.8518
$ws.Range("J132").Value = 3923.111
$ws.Range("K132").Value = 6803.555399999999
$ws.Range("L132").Value = 11769.333
$ws.Range("M132").Value = -4273.555399999999

$ws.Range("H136").Value = 2700.8064
$ws.Range("I136").Value = 1983.1666
$ws.Range("J136").Value = 2873.04
$ws.Range("K136").Value = 5949.4998
$ws.Range("L136").Value = 8619.119999999999
$ws.Range("M136").Value = -3399.4998
$ws.Range("N136").Value = -13719.12

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5236.909
$ws.Range("I134").Value = 4996
$ws.Range("J134").Value = 5879.3335
$ws.Range("K134").Value = 14988
$ws.Range("L134").Value = 17638.0005
$ws.Range("M134").Value = -12453
$ws.Range("N134").Value = -22708.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23439.385
$ws.Range("I31").Value = 1060.7428
$ws.Range("J31").Value = 49547.8
$ws.Range("K31").Value = 1060.7428
$ws.Range("L31").Value = 49547.8
$ws.Range("M31").Value = -765.7428
$ws.Range("N31").Value = -50137.8

$ws.Range("H34").Value = 23439.385
$ws.Range("I34").Value = 1060.7428
$ws.Range("J34").Value = 49547.8
$ws.Range("K34").Value = 1060.7428
$ws.Range("L34").Value = 49547.8
$ws.Range("M34").Value = -858.7428
$ws.Range("N34").Value = -49951.8

$ws.Range("H37").Value = 39900
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 39900
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 39900
$ws.Range("N37").Value = -40114

$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("N45").ClearContents()

$ws.Range("H58").Value = 1117.0209
$ws.Range("I58").Value = 1017.72095
$ws.Range("J58").Value = 1971
$ws.Range("K58").Value = 1017.72095
$ws.Range("L58").Value = 1971
$ws.Range("M58").Value = -814.72095
$ws.Range("N58").Value = -2377

$ws.Range("H92").Value = 30200.666
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 30200.666
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 30200.666
$ws.Range("N92").Value = -35192.666

$ws.Range("H96").Value = 24666.666
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 24666.666
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 24666.666
$ws.Range("N96").Value = -30158.666

$ws.Range("H132").Value = 33336702
$ws.Range("I132").Value = 30306304
$ws.Range("J132").Value = 41670292
$ws.Range("K132").Value = 90918912
$ws.Range("L132").Value = 125010876
$ws.Range("M132").Value = -90916382
$ws.Range("N132").Value = -125015936

$ws.Range("H136").Value = 1117.0209
$ws.Range("I136").Value = 1017.72095
$ws.Range("J136").Value = 1971
$ws.Range("K136").Value = 3053.16285
$ws.Range("L136").Value = 5913
$ws.Range("M136").Value = -503.1628500000002
$ws.Range("N136").Value = -11013

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 855212.3
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 855212.3
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 2565636.9
$ws.Range("N37").Value = -2565860.9

$ws.Range("H51").Value = 2357.1428
$ws.Range("I51").Value = 1100
$ws.Range("J51").Value = 2860
$ws.Range("K51").Value = 3300
$ws.Range("L51").Value = 8580
$ws.Range("M51").Value = -2840
$ws.Range("N51").Value = -9500

$ws.Range("H96").Value = 3875
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 3875
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 11625
$ws.Range("N96").Value = -15743

$ws.Range("H131").Value = 1553.8572
$ws.Range("I131").Value = 2033.3334
$ws.Range("J131").Value = 1536.0988
$ws.Range("K131").Value = 6100.0002
$ws.Range("L131").Value = 4608.2964
$ws.Range("M131").Value = -1060.0002
$ws.Range("N131").Value = -14688.2964

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 50196.79
$ws.Range("I70").Value = 86578.336
$ws.Range("J70").Value = 4241.1577
$ws.Range("K70").Value = 86578.336
$ws.Range("L70").Value = 4241.1577
$ws.Range("M70").Value = -86308.336
$ws.Range("N70").Value = -4781.1577

$ws.Range("H73").Value = 50196.79
$ws.Range("I73").Value = 86578.336
$ws.Range("J73").Value = 4241.1577
$ws.Range("K73").Value = 86578.336
$ws.Range("L73").Value = 4241.1577
$ws.Range("M73").Value = -85642.336
$ws.Range("N73").Value = -6113.1577

$ws.Range("H80").Value = 100104696
$ws.Range("I80").Value = 143005710
$ws.Range("J80").Value = 2300
$ws.Range("K80").Value = 143005710
$ws.Range("L80").Value = 2300
$ws.Range("M80").Value = -143004712
$ws.Range("N80").Value = -4296

$ws.Range("H83").Value = 100104696
$ws.Range("I83").Value = 143005710
$ws.Range("J83").Value = 2300
$ws.Range("K83").Value = 715028550
$ws.Range("L83").Value = 11500
$ws.Range("M83").Value = -715023558
$ws.Range("N83").Value = -21484

$ws.Range("H92").Value = 8998.5
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 8998.5
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 8998.5
$ws.Range("N92").Value = -12742.5

$ws.Range("H94").Value = 18958.8
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 18958.8
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 18958.8
$ws.Range("N94").Value = -20310.8

$ws.Range("H122").Value = 845.6667
$ws.Range("I122").Value = 845.6667
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2537.0001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -87.0001000000002
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 4438.522
$ws.Range("I132").Value = 4260.3335
$ws.Range("J132").Value = 4772.625
$ws.Range("K132").Value = 12781.0005
$ws.Range("L132").Value = 14317.875
$ws.Range("M132").Value = -10251.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 8540
$ws.Range("I35").Value = 1033.3334
$ws.Range("J35").Value = 19800
$ws.Range("K35").Value = 1033.3334
$ws.Range("L35").Value = 19800
$ws.Range("M35").Value = -697.3334
$ws.Range("N35").Value = -20472

$ws.Range("H93").Value = 1565.25
$ws.Range("I93").Value = 1820
$ws.Range("J93").Value = 801
$ws.Range("K93").Value = 1820
$ws.Range("L93").Value = 801
$ws.Range("M93").Value = -572
$ws.Range("N93").Value = -3297

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2354.6365
$ws.Range("I122").Value = 2977.8
$ws.Range("J122").Value = 2171.353
$ws.Range("K122").Value = 8933.400000000001
$ws.Range("L122").Value = 6514.059
$ws.Range("M122").Value = -11414.059
